# Rename the single worksheet from "TDICTADO.RPT" to "data".
# (The defined name "_xlnm._FilterDatabase" references the sheet by name,
# so Excel automatically updates it to point at the renamed sheet.)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TDICTADO.RPT")
$ws.Name = "data"
